$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.007.32"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "1.874.52"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.58"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5079"
$ws.Range("E7").Value = "  -0.79%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3670"
$ws.Range("E8").Value = "  -2.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07213"
$ws.Range("E9").Value = "  +0.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8963"
$ws.Range("E10").Value = "  +0.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.79"
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07525"
$ws.Range("E12").Value = "  -0.61%  "
$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "95.11"
$ws.Range("E13").Value = "  +6.27%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.853.78"
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.254"
$ws.Range("E15").Value = "  -0.97%  "
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008535"
$ws.Range("E17").Value = "  +1.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.27"
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9993"
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").Value = "27.045.81"
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.032"
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").Value = "2.089.03"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.41"
$ws.Range("E23").Value = "  -1.20%  "
$ws.Range("E24").Value = "  -0.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.38"
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.788"
$ws.Range("E26").Value = "  -2.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.91"
$ws.Range("E27").Value = "  -0.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.091"
$ws.Range("E28").Value = "  -0.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.47"
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.725"
$ws.Range("E30").Value = "  +1.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.726"
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09151"
$ws.Range("E32").Value = "  +0.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05124"
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7524"
$ws.Range("E34").Value = "  +3.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.991"
$ws.Range("E35").Value = "  -1.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.162"
$ws.Range("E36").Value = "  +0.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.234"
$ws.Range("E37").Value = "  +6.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.569"
$ws.Range("E38").Value = "  +2.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5695"
$ws.Range("E39").Value = "  +7.11%  "
$ws.Range("E40").Value = "  -1.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.075"
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.633"
$ws.Range("E42").Value = "  +1.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "115.56"
$ws.Range("E43").Value = "  -1.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.595"
$ws.Range("E44").Value = "  +3.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1478"
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4792"
$ws.Range("E46").Value = "  +3.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9990"
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.12"
$ws.Range("E48").Value = "  +1.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.572"
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.91"
$ws.Range("E50").Value = "  +1.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.26"
$ws.Range("E51").Value = "  -0.93%  "
